$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The 'Price' column (D) contains values that look numeric (e.g. '0.998', '596.03')
# but must stay as plain text, exactly like the rest of the sheet. Mark the range as
# Text-formatted first so Excel does not silently convert these into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated crypto price / 1h-volume data scraped by the GitHub Actions job,
# including EnergySwap (row 51) being replaced by Arweave.

$ws.Range("D2").Value = '68.579.14'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '3.785.39'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '596.03'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = '168.32'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '3.784.76'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").Value = '6.51'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '36.83'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = '4.418.96'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '3.775.98'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").Value = '68.501.70'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '18.36'
$ws.Range("E18").Value = '  -3.05%  '
$ws.Range("D19").Value = '7.09'
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '10.88'
$ws.Range("E21").Value = '  +3.02%  '
$ws.Range("D22").Value = '472.21'
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("D23").Value = '0.706'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '85.02'
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("D25").Value = '0.0000146'
$ws.Range("E25").Value = '  -2.93%  '
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").Value = '12.22'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = '10.22'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D30").Value = '3.930.16'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").Value = '2.81'
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("D32").Value = '7.49'
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = '2.26'
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").Value = '30.24'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").Value = '9.27'
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("D36").Value = '0.996'
$ws.Range("D37").Value = '3.739.57'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  -2.40%  '
$ws.Range("D39").Value = '3.53'
$ws.Range("E39").Value = '  -5.74%  '
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = '5.84'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("D44").Value = '0.310'
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '1.96'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '8.62'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '403.29'
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("D49").Value = '45.51'
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").Value = '144.80'
$ws.Range("E50").Value = '  +2.23%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = '40.25'
$ws.Range("E51").Value = '  +3.81%  '

$wb.Save()
